$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RateCompare")
[void]$ws.Activate()

# New column G: "Rct-Name" header + constant reaction id for every data row
$ws.Range("G1").Value = "Rct-Name"
$ws.Range("G2").Value = "EX_glc__D_e"
$ws.Range("G3").Value = "EX_glc__D_e"
$ws.Range("G4").Value = "EX_glc__D_e"
$ws.Range("G5").Value = "EX_glc__D_e"
$ws.Range("G6").Value = "EX_glc__D_e"

# Row 2 (2229v1): updated growth-rate analysis numbers
$ws.Range("C2").Value = 0.03
$ws.Range("D2").Value = 0.01
$ws.Range("E2").Value = 0.49
$ws.Range("F2").Value = 0.14

# Row 5 (LV3_200_v1): Tspan + recomputed numbers (was the "100-200"/0.01/0/0.1/0.02 outlier)
$ws.Range("B5").Value = "0-80"
$ws.Range("C5").Value = 0.07
$ws.Range("D5").Value = 0.01
$ws.Range("E5").Value = 1.06
$ws.Range("F5").Value = 0.38

# Restore the last on-screen selection recorded for this sheet
[void]$ws.Range("K18").Select()
